$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'264.82"
$ws.Range("D3").Value = "'22.46"
$ws.Range("D4").Value = "'6.291"
$ws.Range("D5").Value = "'0.06147"
$ws.Range("D6").Value = "'3.597"
$ws.Range("D7").Value = "'6.666"
$ws.Range("D9").Value = "'0.8290"
$ws.Range("D10").Value = "'0.01357"
$ws.Range("D12").Value = "'0.08232"
$ws.Range("D13").Value = "'0.03429"
$ws.Range("D14").Value = "'0.03130"
$ws.Range("D15").Value = "'0.09246"
$ws.Range("D16").Value = "'3.893"
$ws.Range("D17").Value = "'0.001715"
$ws.Range("D18").Value = "'0.04882"
$ws.Range("D19").Value = "'0.006228"
$ws.Range("D20").Value = "'0.005267"
$ws.Range("D23").Value = "'3.766"
$ws.Range("D24").Value = "'2.289"
$ws.Range("D27").Value = "'0.0002682"
$ws.Range("D40").Value = "'0.04612"
$ws.Range("D41").Value = "'0.006956"
$ws.Range("D43").Value = "'0.003247"
$ws.Range("D44").Value = "'0.01184"
$ws.Range("D45").Value = "'0.00006134"
$ws.Range("D47").Value = "'0.7786"
$ws.Range("D48").Value = "'0.1934"
